$d = $word.ActiveDocument

# Change 1: add line break run after "Hallo {{Vorname}},"
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertBreak(6)
$brRange = $d.Range(18, 19)
$brRange.Style = "Fett"
$brRange.Font.Bold = 0
$brRange.Font.BoldBi = 0
Write-Output "change1 done"

# Change 2: merge "Ergebnis der " + "UStVA" (remove proofErr wrapping)
$r2a = $d.Content
$found2a = $r2a.Find.Execute("Ergebnis der UStVA", $true, $false, $false, $false, $false, $true, 1, $false, "Ergebnis der UStVAQQQ", 2)
Write-Output "found2a=$found2a"
$r2b = $d.Content
$found2b = $r2b.Find.Execute("UStVAQQQ", $true, $false, $false, $false, $false, $true, 1, $false, "UStVA", 2)
Write-Output "found2b=$found2b"

# Change 3: merge "Addison " + "OneClick"
$r3a = $d.Content
$found3a = $r3a.Find.Execute("Addison OneClick", $true, $false, $false, $false, $false, $true, 1, $false, "Addison OneClickQQQ", 2)
Write-Output "found3a=$found3a"
$r3b = $d.Content
$found3b = $r3b.Find.Execute("OneClickQQQ", $true, $false, $false, $false, $false, $true, 1, $false, "OneClick", 2)
Write-Output "found3b=$found3b"

# Change 4: insert "für einen Beratungstermin " before "bei uns"
$r4 = $d.Content
$found4 = $r4.Find.Execute("melde Dich gerne ")
Write-Output "found4=$found4"
$r4.Collapse(0)
$r4.InsertAfter("für einen Beratungstermin ")
Write-Output "change4 done"
